# Crop Rotation Sample Template - "begin teting with multiple year data"
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Farmland sheet: mark rows 4-11 as Included ("Y"), rename "Row 10" -> "Row A"
# ---------------------------------------------------------------------------
$wsFarm = $wb.Worksheets.Item("Farmland")
$wsFarm.Activate()

$wsFarm.Range("D4").Value = "Y"
$wsFarm.Range("D5").Value = "Y"
$wsFarm.Range("D6").Value = "Y"
$wsFarm.Range("D7").Value = "Y"
$wsFarm.Range("D8").Value = "Y"
$wsFarm.Range("D9").Value = "Y"
$wsFarm.Range("D10").Value = "Y"
$wsFarm.Range("D11").Value = "Y"

# Renaming the shared string updates every other reference to "Row 10"
# across the workbook (there are none besides this cell).
$wsFarm.Range("A11").Value = "Row A"

$wsFarm.Range("A12").Select()

Write-Host "Farmland sheet updated"

# ---------------------------------------------------------------------------
# 2) Crops sheet: mark several crops "Plant?" = Y, move the view
# ---------------------------------------------------------------------------
$wsCrops = $wb.Worksheets.Item("Crops")
$wsCrops.Activate()

$wsCrops.Range("X7").Value = "Y"
$wsCrops.Range("X8").Value = "Y"
$wsCrops.Range("X9").Value = "Y"
$wsCrops.Range("X13").Value = "Y"
$wsCrops.Range("X37").Value = "Y"
$wsCrops.Range("X38").Value = "Y"
$wsCrops.Range("X39").Value = "Y"
$wsCrops.Range("X44").Value = "Y"

$wsCrops.Range("B3").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsCrops.Range("D9").Select()

Write-Host "Crops sheet updated"
